$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated cryptos list (prices and 1h volume change %) ---

# Cells below contain numeric-looking text (e.g. "0.999", "37.10") that must
# stay plain text (matching the sheet's existing inline-string data), so we
# force a Text number format before writing the value, just like Excel does
# when a user pre-formats a column as Text before typing numeric-looking data.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D9",
    "D12",
    "D13",
    "D14",
    "D17",
    "D19",
    "D21",
    "D22",
    "D25",
    "D27",
    "D28",
    "D32",
    "D34",
    "D35",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D45",
    "D46",
    "D47",
    "D48",
    "D50",
    "D51",
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$values = @{
    "D4" = "0.999"
    "D5" = "600.26"
    "D6" = "168.97"
    "D9" = "0.531"
    "D12" = "0.462"
    "D13" = "0.0000275"
    "D14" = "37.10"
    "D17" = "18.75"
    "D19" = "7.47"
    "D21" = "10.86"
    "D22" = "469.32"
    "D25" = "83.74"
    "D27" = "12.21"
    "D28" = "10.36"
    "D32" = "7.75"
    "D34" = "30.72"
    "D35" = "9.32"
    "D39" = "5.98"
    "D40" = "1.01"
    "D41" = "0.139"
    "D42" = "0.999"
    "D43" = "0.319"
    "D45" = "1.98"
    "D46" = "8.77"
    "D47" = "407.92"
    "D48" = "46.47"
    "D50" = "142.92"
    "D51" = "0.0360"
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

# --- Remaining cells (coin name/link swap and all other price/volume text) ---
$ws.Range("D2").Value = "67.874.17"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "3.819.25"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "3.818.30"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("E13").Value = "  +7.50%  "
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "4.458.36"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").Value = "3.804.25"
$ws.Range("E16").Value = "  -2.48%  "
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").Value = "67.870.13"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("E24").Value = "  -8.92%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("D31").Value = "3.969.14"
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("E34").Value = "  -2.43%  "
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("D36").Value = "3.783.26"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("E37").Value = "  +5.32%  "
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("E47").Value = "  -4.41%  "
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("E49").Value = "  -6.80%  "
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("E51").Value = "  +0.16%  "
